$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: re-run entry whose font color is normalized to black (matches
# the other "weighted" rows), which also bumps the row height back to 19.5 ---
$ws.Rows.Item(29).RowHeight = 19.5

$ws.Range("K29").Font.Color = 0
$ws.Range("N29").Font.Color = 0
$ws.Range("O29").Font.Color = 0

# --- Row 30: new "subrun fts" 10%-subset run, added below row 29 ---
$ws.Range("A30").Value = "ukb51139_subset.csv"
$ws.Range("B30").Value = "2801 x 462"
$ws.Range("C30").Value = "subrun fts"
$ws.Range("D30").Value = "no events"
$ws.Range("E30").Value = "> 140/80"
$ws.Range("F30").Value = "zscore"
$ws.Range("G30").Value = "median"
$ws.Range("H30").Value = "none"
$ws.Range("I30").Value = 50
$ws.Range("K30").Value = 114
$ws.Range("L30").Value = "-230.9 & -42.0"
$ws.Range("M30").Value = "45.8 & 40.4"
$ws.Range("N30").Value = 19
$ws.Range("O30").Value = 68.9
